$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.935.29'
$ws.Range('E2').Value = '  +0.68%  '
$ws.Range('D3').Value = '3.136.05'
$ws.Range('E3').Value = '  +2.21%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '527.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.83'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.51%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.136.60'
$ws.Range('E8').Value = '  +2.24%  '
$ws.Range('E9').Value = '  +0.26%  '
$ws.Range('E10').Value = '  +0.80%  '
$ws.Range('E11').Value = '  +2.37%  '
$ws.Range('E12').Value = '  +3.49%  '
$ws.Range('D13').Value = '3.677.82'
$ws.Range('E13').Value = '  +2.12%  '
$ws.Range('E14').Value = '  +1.60%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.42'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000165'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.60%  '
$ws.Range('D17').Value = '58.043.55'
$ws.Range('E17').Value = '  +0.74%  '
$ws.Range('D18').Value = '3.139.54'
$ws.Range('E18').Value = '  +2.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.16'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.98'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.15'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.61%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '338.07'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.86%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.514'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.97'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.87%  '
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('E27').Value = '  +0.25%  '
$ws.Range('D28').Value = '0.0₃0936'
$ws.Range('E28').Value = '  +3.93%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.64'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.32%  '
$ws.Range('E30').Value = '  +0.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.27'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.89'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.75%  '
$ws.Range('E33').Value = '  +3.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '21.08'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.53%  '
$ws.Range('B35').Value = 'Monero'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '155.56'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.66%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.70'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.75%  '
$ws.Range('E37').Value = '  +3.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '27.24'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.73%  '
$ws.Range('E39').Value = '  +4.35%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0669'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.44%  '
$ws.Range('D41').Value = '3.180.58'
$ws.Range('E41').Value = '  +2.13%  '
$ws.Range('E42').Value = '  +5.98%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.54'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +13.17%  '
$ws.Range('E44').Value = '  +0.65%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '37.09'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.59%  '
$ws.Range('E46').Value = '  -0.10%  '
$ws.Range('D47').Value = '2.304.11'
$ws.Range('E47').Value = '  +1.91%  '
$ws.Range('E48').Value = '  +1.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.999'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.76%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '21.14'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.03'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.84%  '

Write-Output "Applied changes"
